$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Search")

# Insert a new row above the existing data row (shifts "summer dress" from A1 to A2)
$ws.Rows.Item(1).Insert()

# New header/search-data cell with a yellow highlight fill
$ws.Range("A1").Value = "search data"
$ws.Range("A1").Interior.Color = 65535

# Selection moves to the (now) second row
$ws.Range("A2").Select()

# Page setup: A4, portrait
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
